# Updates the cryptocurrency price/volume table (and the two swapped
# rows for USDC/XRP) to match the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# U+2083 SUBSCRIPT THREE used in the Dogecoin-style "leading zero count"
# price notation for D20 (e.g. 0.0[3]0886).
$subThree = [char]0x2083

$updates = @(
    @{ Cell = "D2"; Value = "42.132.95" },
    @{ Cell = "E2"; Value = "  -2.06%  " },
    @{ Cell = "D3"; Value = "2.268.35" },
    @{ Cell = "E3"; Value = "  -3.05%  " },
    @{ Cell = "E4"; Value = "  +0.02%  " },
    @{ Cell = "D5"; Value = "297.80" },
    @{ Cell = "E5"; Value = "  -2.88%  " },
    @{ Cell = "D6"; Value = "94.32" },
    @{ Cell = "E6"; Value = "  -6.72%  " },
    @{ Cell = "B7"; Value = "USDC" },
    @{ Cell = "C7"; Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc" },
    @{ Cell = "D7"; Value = "1.00" },
    @{ Cell = "E7"; Value = "  +0.09%  " },
    @{ Cell = "B8"; Value = "XRP" },
    @{ Cell = "C8"; Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp" },
    @{ Cell = "D8"; Value = "0.493" },
    @{ Cell = "E8"; Value = "  -3.97%  " },
    @{ Cell = "D9"; Value = "0.488" },
    @{ Cell = "E9"; Value = "  -4.62%  " },
    @{ Cell = "D10"; Value = "33.00" },
    @{ Cell = "E10"; Value = "  -5.59%  " },
    @{ Cell = "D11"; Value = "0.0786" },
    @{ Cell = "E11"; Value = "  -1.84%  " },
    @{ Cell = "D12"; Value = "48.28" },
    @{ Cell = "E12"; Value = "  -7.75%  " },
    @{ Cell = "E13"; Value = "  -0.27%  " },
    @{ Cell = "D14"; Value = "6.62" },
    @{ Cell = "E14"; Value = "  -3.39%  " },
    @{ Cell = "D15"; Value = "15.60" },
    @{ Cell = "E15"; Value = "  -1.78%  " },
    @{ Cell = "D16"; Value = "2.620.95" },
    @{ Cell = "E16"; Value = "  -3.08%  " },
    @{ Cell = "D17"; Value = "2.266.56" },
    @{ Cell = "E17"; Value = "  -1.88%  " },
    @{ Cell = "D18"; Value = "0.771" },
    @{ Cell = "E18"; Value = "  -5.18%  " },
    @{ Cell = "D19"; Value = "42.129.58" },
    @{ Cell = "E19"; Value = "  -1.91%  " },
    @{ Cell = "D20"; Value = ("0.0{0}0886" -f $subThree) },
    @{ Cell = "E20"; Value = "  -2.92%  " },
    @{ Cell = "D21"; Value = "11.28" },
    @{ Cell = "E21"; Value = "  -4.68%  " },
    @{ Cell = "D22"; Value = "5.95" },
    @{ Cell = "E22"; Value = "  -4.57%  " },
    @{ Cell = "D23"; Value = "66.52" },
    @{ Cell = "E23"; Value = "  -1.98%  " },
    @{ Cell = "D24"; Value = "232.40" },
    @{ Cell = "E24"; Value = "  -2.02%  " },
    @{ Cell = "E25"; Value = "  -4.79%  " },
    @{ Cell = "D26"; Value = "0.999" },
    @{ Cell = "E26"; Value = "  -0.05%  " },
    @{ Cell = "D27"; Value = "2.44" },
    @{ Cell = "E27"; Value = "  -4.53%  " },
    @{ Cell = "D28"; Value = "23.72" },
    @{ Cell = "E28"; Value = "  -7.30%  " },
    @{ Cell = "D29"; Value = "2.29" },
    @{ Cell = "E29"; Value = "  -1.29%  " },
    @{ Cell = "D30"; Value = "166.44" },
    @{ Cell = "E30"; Value = "  +1.85%  " },
    @{ Cell = "D31"; Value = "33.51" },
    @{ Cell = "E31"; Value = "  -4.72%  " },
    @{ Cell = "D32"; Value = "9.00" },
    @{ Cell = "E32"; Value = "  -4.25%  " },
    @{ Cell = "D33"; Value = "1.00" },
    @{ Cell = "E33"; Value = "  +0.02%  " },
    @{ Cell = "D34"; Value = "4.90" },
    @{ Cell = "E34"; Value = "  -4.59%  " },
    @{ Cell = "E35"; Value = "  -3.02%  " },
    @{ Cell = "E36"; Value = "  -5.56%  " },
    @{ Cell = "D37"; Value = "0.0687" },
    @{ Cell = "E37"; Value = "  -5.80%  " },
    @{ Cell = "D38"; Value = "16.00" },
    @{ Cell = "E38"; Value = "  -9.54%  " },
    @{ Cell = "D39"; Value = "2.76" },
    @{ Cell = "E39"; Value = "  -5.56%  " },
    @{ Cell = "D40"; Value = "0.0986" },
    @{ Cell = "E40"; Value = "  -3.72%  " },
    @{ Cell = "E41"; Value = "  -4.16%  " },
    @{ Cell = "D42"; Value = "1.70" },
    @{ Cell = "E42"; Value = "  -8.57%  " },
    @{ Cell = "E43"; Value = "  -2.03%  " },
    @{ Cell = "D44"; Value = "1.956.10" },
    @{ Cell = "E44"; Value = "  -2.58%  " },
    @{ Cell = "D45"; Value = "0.0276" },
    @{ Cell = "E45"; Value = "  -3.57%  " },
    @{ Cell = "D46"; Value = "17.28" },
    @{ Cell = "E46"; Value = "  -7.60%  " },
    @{ Cell = "D47"; Value = "9.52" },
    @{ Cell = "E47"; Value = "  -6.70%  " },
    @{ Cell = "D48"; Value = "2.77" },
    @{ Cell = "E48"; Value = "  -5.95%  " },
    @{ Cell = "D49"; Value = "2.493.26" },
    @{ Cell = "E49"; Value = "  -2.46%  " },
    @{ Cell = "E50"; Value = "  -3.38%  " },
    @{ Cell = "D51"; Value = "51.68" },
    @{ Cell = "E51"; Value = "  -7.87%  " }
)

foreach ($u in $updates) {
    # Force text formatting first so numeric-looking strings (prices
    # like "297.80" or "1.00") keep their exact text, matching the
    # original inline-string cells instead of being coerced to numbers.
    $ws.Range($u.Cell).NumberFormat = "@"
    $ws.Range($u.Cell).Value = $u.Value
    # Reset the format back to General/default so we do not leave a
    # stray style behind on cells that originally had none.
    $ws.Range($u.Cell).ClearFormats()
}
